# Update the marksheet totals on the "quiz" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct-answer mark value changed from 3 to 5
$ws.Range("B11").Value = 5

# Total row: total correct marks changed from 66 to 110 (22 right * 5)
$ws.Range("B12").Value = 110

# Total row: obtained/maximum marks label changed from "63/84" to "110/140"
$ws.Range("E12").Value = "110/140"
